$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 374792.72
$ws.Range("I98").Value = 4404.7617
$ws.Range("J98").Value = 1671150.5
$ws.Range("K98").Value = 4404.7617
$ws.Range("L98").Value = 1671150.5
$ws.Range("M98").Value = -2906.7617
$ws.Range("N98").Value = -1674146.5
$ws.Range("H122").Value = 374792.72
$ws.Range("I122").Value = 4404.7617
$ws.Range("J122").Value = 1671150.5
$ws.Range("K122").Value = 13214.2851
$ws.Range("L122").Value = 5013451.5
$ws.Range("M122").Value = -10764.2851
$ws.Range("N122").Value = -5018351.5
$ws.Range("H135").Value = 2070.9607
$ws.Range("I135").Value = 962.5682
$ws.Range("K135").Value = 8663.113800000001
$ws.Range("M135").Value = -6128.113800000001
$ws.Range("H137").Value = 927.4706
$ws.Range("I137").Value = 820.6111
$ws.Range("J137").Value = 1047.6875
$ws.Range("K137").Value = 2461.8333
$ws.Range("L137").Value = 3143.0625
$ws.Range("M137").Value = 88.16670000000022
$ws.Range("N137").Value = -8243.0625
$ws.Range("H138").Value = 3641.6936
$ws.Range("I138").Value = 2135.1333
$ws.Range("J138").Value = 5054.0938
$ws.Range("K138").Value = 6405.3999
$ws.Range("L138").Value = 15162.2814
$ws.Range("M138").Value = -1265.3999
$ws.Range("N138").Value = -25442.2814
$ws.Range("H141").Value = 2072.9429
$ws.Range("I141").Value = 2045.6765
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 6137.029500000001
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -957.0295000000006
$ws.Range("N141").Value = -19360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 418055.22
$ws.Range("I32").Value = 3172
$ws.Range("J32").Value = 3581539.8
$ws.Range("K32").Value = 3172
$ws.Range("L32").Value = 3581539.8
$ws.Range("M32").Value = -2885
$ws.Range("N32").Value = -3582113.8
$ws.Range("H61").Value = 1115.3334
$ws.Range("I61").Value = 861.2414
$ws.Range("J61").Value = 1682.1538
$ws.Range("K61").Value = 861.2414
$ws.Range("L61").Value = 1682.1538
$ws.Range("M61").Value = -649.2414
$ws.Range("N61").Value = -2106.1538
$ws.Range("H97").Value = 1369.4242
$ws.Range("I97").Value = 834
$ws.Range("K97").Value = 834
$ws.Range("M97").Value = -338
$ws.Range("H136").Value = 1115.3334
$ws.Range("I136").Value = 861.2414
$ws.Range("J136").Value = 1682.1538
$ws.Range("K136").Value = 2583.7242
$ws.Range("L136").Value = 5046.4614
$ws.Range("M136").Value = -33.72420000000011
$ws.Range("N136").Value = -10146.4614

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6852086
$ws.Range("I31").Value = 8773513
$ws.Range("J31").Value = 7001.125
$ws.Range("K31").Value = 8773513
$ws.Range("L31").Value = 7001.125
$ws.Range("M31").Value = -8773218
$ws.Range("N31").Value = -7591.125
$ws.Range("H34").Value = 6852086
$ws.Range("I34").Value = 8773513
$ws.Range("J34").Value = 7001.125
$ws.Range("K34").Value = 8773513
$ws.Range("L34").Value = 7001.125
$ws.Range("M34").Value = -8773311
$ws.Range("N34").Value = -7405.125
$ws.Range("H134").Value = 1357.234
$ws.Range("I134").Value = 1325.7106
$ws.Range("J134").Value = 1490.3334
$ws.Range("K134").Value = 3977.1318
$ws.Range("L134").Value = 4471.0002
$ws.Range("M134").Value = -1442.1318
$ws.Range("N134").Value = -9541.0002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1245.1538
$ws.Range("I132").Value = 1190.25
$ws.Range("J132").Value = 1302.9474
$ws.Range("K132").Value = 10712.25
$ws.Range("L132").Value = 11726.5266
$ws.Range("M132").Value = -8182.25
$ws.Range("N132").Value = -16786.5266

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 33.555557
$ws.Range("I2").Value = 28.75
$ws.Range("J2").Value = 37.4
$ws.Range("K2").Value = 28.75
$ws.Range("L2").Value = 37.4
$ws.Range("M2").Value = 84.25
$ws.Range("N2").Value = -263.4
$ws.Range("H43").Value = 6342301.5
$ws.Range("I43").Value = 19000150
$ws.Range("J43").Value = 13377
$ws.Range("K43").Value = 19000150
$ws.Range("L43").Value = 13377
$ws.Range("M43").Value = -18999999
$ws.Range("N43").Value = -13679
$ws.Range("H46").Value = 13795
$ws.Range("I46").Value = 10934.429
$ws.Range("J46").Value = 17799.8
$ws.Range("K46").Value = 10934.429
$ws.Range("L46").Value = 17799.8
$ws.Range("M46").Value = -10778.429
$ws.Range("N46").Value = -18111.8
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640
$ws.Range("H70").Value = 14071771
$ws.Range("I70").Value = 35168920
$ws.Range("J70").Value = 7005.5415
$ws.Range("K70").Value = 35168920
$ws.Range("L70").Value = 7005.5415
$ws.Range("M70").Value = -35168650
$ws.Range("N70").Value = -7545.5415
$ws.Range("H73").Value = 14071771
$ws.Range("I73").Value = 35168920
$ws.Range("J73").Value = 7005.5415
$ws.Range("K73").Value = 35168920
$ws.Range("L73").Value = 7005.5415
$ws.Range("M73").Value = -35167984
$ws.Range("N73").Value = -8877.541499999999
$ws.Range("H80").Value = 2864.6487
$ws.Range("I80").Value = 2532.6667
$ws.Range("J80").Value = 3477.5386
$ws.Range("K80").Value = 2532.6667
$ws.Range("L80").Value = 3477.5386
$ws.Range("M80").Value = -1534.6667
$ws.Range("N80").Value = -5473.5386
$ws.Range("H83").Value = 2864.6487
$ws.Range("I83").Value = 2532.6667
$ws.Range("J83").Value = 3477.5386
$ws.Range("K83").Value = 12663.3335
$ws.Range("L83").Value = 17387.693
$ws.Range("M83").Value = -7671.333500000001
$ws.Range("N83").Value = -27371.693

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2414.328
$ws.Range("I136").Value = 1846.75
$ws.Range("J136").Value = 4510
$ws.Range("K136").Value = 5540.25
$ws.Range("L136").Value = 13530
$ws.Range("M136").Value = -2990.25
$ws.Range("N136").Value = -18630

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4241.6665
$ws.Range("I96").Value = 1750
$ws.Range("K96").Value = 1750
$ws.Range("M96").Value = -377
